# Updated cryptocurrency Price/Volume(1h) data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells are authored as literal text (e.g. "34.748.68",
# "11.00", "0.1000") - not real numbers. Force Text format first so Excel
# COM does not renormalize/round them when the .Value is assigned.
$priceCells = @(
    'D2', 'D3', 'D5', 'D6', 'D8', 'D10', 'D11', 'D12',
    'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20',
    'D21', 'D22', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29',
    'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38',
    'D39', 'D40', 'D42', 'D43', 'D44', 'D46', 'D47', 'D48',
    'D49', 'D51'
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range('D2').Value = '34.748.68'
$ws.Range('E2').Value = '  +0.68%  '
# Row 3
$ws.Range('D3').Value = '1.807.73'
$ws.Range('E3').Value = '  +0.31%  '
# Row 4
$ws.Range('E4').Value = '  +0.67%  '
# Row 5
$ws.Range('D5').Value = '230.69'
$ws.Range('E5').Value = '  +2.83%  '
# Row 6
$ws.Range('D6').Value = '0.603'
$ws.Range('E6').Value = '  +0.47%  '
# Row 7
$ws.Range('E7').Value = '  +0.60%  '
# Row 8
$ws.Range('D8').Value = '40.47'
$ws.Range('E8').Value = '  -1.66%  '
# Row 9
$ws.Range('E9').Value = '  +4.49%  '
# Row 10
$ws.Range('D10').Value = '0.0678'
$ws.Range('E10').Value = '  +1.78%  '
# Row 11
$ws.Range('D11').Value = '0.1000'
$ws.Range('E11').Value = '  +0.23%  '
# Row 12
$ws.Range('D12').Value = '2.077.61'
$ws.Range('E12').Value = '  +0.77%  '
# Row 13
$ws.Range('D13').Value = '1.813.38'
$ws.Range('E13').Value = '  +0.83%  '
# Row 14
$ws.Range('D14').Value = '11.00'
$ws.Range('E14').Value = '  +0.47%  '
# Row 15
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = '0.652'
$ws.Range('E15').Value = '  +3.85%  '
# Row 16
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').Value = '4.60'
$ws.Range('E16').Value = '  +4.54%  '
# Row 17
$ws.Range('D17').Value = '34.775.74'
$ws.Range('E17').Value = '  +0.84%  '
# Row 18
$ws.Range('D18').Value = '68.72'
$ws.Range('E18').Value = '  +2.28%  '
# Row 19
$ws.Range('D19').Value = '0.0₃0779'
$ws.Range('E19').Value = '  +1.56%  '
# Row 20
$ws.Range('D20').Value = '235.88'
$ws.Range('E20').Value = '  -1.72%  '
# Row 21
$ws.Range('D21').Value = '11.66'
$ws.Range('E21').Value = '  +4.89%  '
# Row 22
$ws.Range('D22').Value = '4.60'
$ws.Range('E22').Value = '  +9.03%  '
# Row 23
$ws.Range('E23').Value = '  +0.45%  '
# Row 24
$ws.Range('D24').Value = '2.23'
$ws.Range('E24').Value = '  +3.32%  '
# Row 25
$ws.Range('D25').Value = '172.21'
$ws.Range('E25').Value = '  +0.08%  '
# Row 26
$ws.Range('D26').Value = '7.69'
$ws.Range('E26').Value = '  -0.13%  '
# Row 27
$ws.Range('D27').Value = '17.22'
$ws.Range('E27').Value = '  -0.98%  '
# Row 28
$ws.Range('D28').Value = '0.119'
$ws.Range('E28').Value = '  -0.87%  '
# Row 29
$ws.Range('D29').Value = '1.58'
$ws.Range('E29').Value = '  +29.44%  '
# Row 30
$ws.Range('E30').Value = '  +0.47%  '
# Row 31
$ws.Range('D31').Value = '3.341.15'
$ws.Range('E31').Value = '  +37.51%  '
# Row 32
$ws.Range('D32').Value = '0.0541'
$ws.Range('E32').Value = '  +5.76%  '
# Row 33
$ws.Range('D33').Value = '3.85'
$ws.Range('E33').Value = '  +1.85%  '
# Row 34
$ws.Range('D34').Value = '3.95'
$ws.Range('E34').Value = '  +2.52%  '
# Row 35
$ws.Range('D35').Value = '1.75'
$ws.Range('E35').Value = '  -2.37%  '
# Row 36
$ws.Range('D36').Value = '92.29'
$ws.Range('E36').Value = '  +7.26%  '
# Row 37
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '1.11'
$ws.Range('E37').Value = '  +5.01%  '
# Row 38
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '0.669'
$ws.Range('E38').Value = '  +3.42%  '
# Row 39
$ws.Range('D39').Value = '1.300.84'
$ws.Range('E39').Value = '  -1.65%  '
# Row 40
$ws.Range('D40').Value = '0.0190'
$ws.Range('E40').Value = '  +1.57%  '
# Row 41
$ws.Range('E41').Value = '  +3.43%  '
# Row 42
$ws.Range('D42').Value = '0.973'
$ws.Range('E42').Value = '  +3.89%  '
# Row 43
$ws.Range('D43').Value = '14.63'
$ws.Range('E43').Value = '  -0.52%  '
# Row 44
$ws.Range('D44').Value = '2.30'
$ws.Range('E44').Value = '  -2.07%  '
# Row 45
$ws.Range('E45').Value = '  +0.61%  '
# Row 46
$ws.Range('D46').Value = '2.74'
$ws.Range('E46').Value = '  -1.65%  '
# Row 47
$ws.Range('D47').Value = '6.15'
$ws.Range('E47').Value = '  +6.03%  '
# Row 48
$ws.Range('D48').Value = '0.0510'
$ws.Range('E48').Value = '  -1.75%  '
# Row 49
$ws.Range('D49').Value = '1.991.56'
$ws.Range('E49').Value = '  +1.53%  '
# Row 50
$ws.Range('E50').Value = '  +0.57%  '
# Row 51
$ws.Range('D51').Value = '99.38'
$ws.Range('E51').Value = '  -1.05%  '
